# Fruta / hortaliza, semanal
# Insert 4 new weekly price rows (Valencia orange, Vega Modelo de Temuco)
# above the existing row 570, pushing the rest of the table down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 570 (old rows 570-623 become 574-627)
$ws.Rows("570:573").Insert()

# Constant columns shared by every record in this sheet subset
$mercadoId = 10
$mercado   = "Vega Modelo de Temuco"
$region    = "La Araucanía"
$codreg    = 9
$tipo      = "Fruta"
$prodId    = 100102
$producto  = "Cítricos"
$catId     = 100102005
$categoria = "Naranja"
$origenRegOHiggins = "Región de O'Higgins"

# New rows: Fecha, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm,
#           UnidadComercializacion, Origen, PrecioKg, KgUnidad
$newRows = @(
    @{ Row=570; Fecha=44578; Variedad="Valencia"; Calidad="Especial"; Volumen=65;  PMin=18000; PMax=18000; PProm=18000; Unidad="$/caja 15 kilos granel";     Origen=$origenRegOHiggins; PKg=1200; KgU=15 },
    @{ Row=571; Fecha=44578; Variedad="Valencia"; Calidad="Primera";  Volumen=135; PMin=15000; PMax=15000; PProm=15000; Unidad="$/bandeja 15 kilos granel"; Origen=$origenRegOHiggins; PKg=1000; KgU=15 },
    @{ Row=572; Fecha=44578; Variedad="Valencia"; Calidad="Segunda";  Volumen=125; PMin=12000; PMax=12000; PProm=12000; Unidad="$/bandeja 15 kilos granel"; Origen=$origenRegOHiggins; PKg=800;  KgU=15 },
    @{ Row=573; Fecha=44578; Variedad="Valencia"; Calidad="Tercera";  Volumen=110; PMin=9000;  PMax=9000;  PProm=9000;  Unidad="$/bandeja 15 kilos granel"; Origen=$origenRegOHiggins; PKg=600;  KgU=15 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $prodId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $catId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.KgU
}
